$d = $word.ActiveDocument
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:tbl><w:tblPr><w:tblW w:w="9746" w:type="dxa"/><w:tblInd w:w="34" w:type="dxa"/><w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tblBorders><w:tblCellMar><w:left w:w="70" w:type="dxa"/><w:right w:w="70" w:type="dxa"/></w:tblCellMar><w:tblLook w:val="0000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0"/></w:tblPr><w:tblGrid><w:gridCol w:w="4068"/><w:gridCol w:w="1016"/><w:gridCol w:w="2980"/><w:gridCol w:w="876"/><w:gridCol w:w="806"/></w:tblGrid><w:tr w:rsidR="00BA1882" w14:paraId="498B11C7" w14:textId="77777777" w:rsidTr="00C8027A"><w:trPr><w:trHeight w:val="211"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="9746" w:type="dxa"/><w:gridSpan w:val="5"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/></w:tcPr><w:p w14:paraId="30BC796D" w14:textId="77777777" w:rsidR="00BA1882" w:rsidRPr="00BA1882" w:rsidRDefault="00BA1882" w:rsidP="00FB5666"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Solicitud de cambio</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00CA29AE" w14:paraId="56573E63" w14:textId="49E94D4D" w:rsidTr="00CA29AE"><w:trPr><w:trHeight w:val="211"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="4068" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/></w:tcPr><w:p w14:paraId="58AF40D6" w14:textId="1E1A350A" w:rsidR="00CA29AE" w:rsidRDefault="00CA29AE" w:rsidP="00CA29AE"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Equipo al que va dirigida la solicitud:</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3996" w:type="dxa"/><w:gridSpan w:val="2"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/></w:tcPr><w:p w14:paraId="16476A4D" w14:textId="77777777" w:rsidR="00CA29AE" w:rsidRDefault="00CA29AE" w:rsidP="00CA29AE"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="876" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>#ID:</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="806" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:tc></w:tr><w:tr w:rsidR="00BA1882" w14:paraId="6B0DDA09" w14:textId="77777777" w:rsidTr="00C8027A"><w:trPr><w:trHeight w:val="251"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="9746" w:type="dxa"/><w:gridSpan w:val="5"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/></w:tcPr><w:p w14:paraId="7735EE8F" w14:textId="77777777" w:rsidR="00BA1882" w:rsidRPr="00224C83" w:rsidRDefault="00224C83" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="00224C83"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Propósito del cambio:</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00A75C88" w14:paraId="7A24A354" w14:textId="77777777" w:rsidTr="00C8027A"><w:trPr><w:trHeight w:val="865"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="9746" w:type="dxa"/><w:gridSpan w:val="5"/></w:tcPr><w:p w14:paraId="0D726F4D" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRPr="00224C83" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr w:rsidR="007307BE" w14:paraId="7BAB98AE" w14:textId="0A956B00" w:rsidTr="006A7529"><w:trPr><w:trHeight w:val="671"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5084" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p w14:paraId="6A9ABCDB" w14:textId="77777777" w:rsidR="007307BE" w:rsidRDefault="007307BE" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251668480" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="5DDE0502" wp14:editId="29731637"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>1950085</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>267335</wp:posOffset></wp:positionV><wp:extent cx="297180" cy="251460"/><wp:effectExtent l="0" t="0" r="26670" b="15240"/><wp:wrapNone/><wp:docPr id="5" name="Rectángulo 5"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="297180" cy="251460"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:solidFill><a:schemeClr val="tx1"/></a:solidFill></a:ln></wps:spPr><wps:style><a:lnRef idx="2"><a:schemeClr val="accent1"><a:shade val="50000"/></a:schemeClr></a:lnRef><a:fillRef idx="1"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect w14:anchorId="5319FB74" id="Rectángulo 5" o:spid="_x0000_s1026" style="position:absolute;margin-left:153.55pt;margin-top:21.05pt;width:23.4pt;height:19.8pt;z-index:251668480;visibility:visible;mso-wrap-style:square;mso-height-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-height-percent:0;mso-height-relative:margin;v-text-anchor:middle" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQBaik2voAIAAI4FAAAOAAAAZHJzL2Uyb0RvYy54bWysVM1u2zAMvg/YOwi6r46DuG2COkXQosOA&#10;oi36g55VWYoNyKImKXGyt9mz7MVGSbYTdMUOw3JQRJP8KH78ubjctYpshXUN6JLmJxNKhOZQNXpd&#10;0pfnmy/nlDjPdMUUaFHSvXD0cvn500VnFmIKNahKWIIg2i06U9Lae7PIMsdr0TJ3AkZoVEqwLfMo&#10;2nVWWdYhequy6WRymnVgK2OBC+fw63VS0mXEl1Jwfy+lE56okuLbfDxtPN/CmS0v2GJtmakb3j+D&#10;/cMrWtZoDDpCXTPPyMY2f0C1DbfgQPoTDm0GUjZcxBwwm3zyLpunmhkRc0FynBlpcv8Plt9tHyxp&#10;qpIWlGjWYokekbRfP/V6o4AUgaDOuAXaPZkH20sOryHbnbRt+Mc8yC6Suh9JFTtPOH6czs/yc6Se&#10;o2pa5LPTSHp2cDbW+a8CWhIuJbUYPlLJtrfOY0A0HUxCLA03jVKxbkqTrqTzYlpEBweqqYIymMUO&#10;ElfKki3D2vtdHlJBrCMrlJTGjyHBlFK8+b0SAULpRyGRm5BEChC68oDJOBfa50lVs0qkUMUEf0Ow&#10;wSOGjoABWeIjR+weYLBMIAN2enNvH1xFbOrRefK3hyXn0SNGBu1H57bRYD8CUJhVHznZDyQlagJL&#10;b1DtsXMspJFyht80WL9b5vwDszhDWHLcC/4eD6kA6wT9jZIa7I+Pvgd7bG3UUtLhTJbUfd8wKyhR&#10;3zQ2/TyfzcIQR2FWnE1RsMeat2ON3rRXgKXPcQMZHq/B3qvhKi20r7g+ViEqqpjmGLuk3NtBuPJp&#10;V+AC4mK1imY4uIb5W/1keAAPrIb+fN69Mmv6JvbY/XcwzC9bvOvlZBs8Naw2HmQTG/3Aa883Dn1s&#10;nH5Bha1yLEerwxpd/gYAAP//AwBQSwMEFAAGAAgAAAAhAHIrnUrdAAAACQEAAA8AAABkcnMvZG93&#10;bnJldi54bWxMj8tOwzAQRfdI/IM1SGxQ6zygDSFOhZBYBonCB7jxEEf1q7HThr9nWMFqNJqjO+c2&#10;u8UadsYpjt4JyNcZMHS9V6MbBHx+vK4qYDFJp6TxDgV8Y4Rde33VyFr5i3vH8z4NjEJcrKUAnVKo&#10;OY+9Rivj2gd0dPvyk5WJ1mngapIXCreGF1m24VaOjj5oGfBFY3/cz1bAMlenUzcfrcayM3dFCm9d&#10;CELc3izPT8ASLukPhl99UoeWnA5+dioyI6DMtjmhAu4LmgSUD+UjsIOAKt8Cbxv+v0H7AwAA//8D&#10;AFBLAQItABQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAAAAAAAAAAAAAAAAAAABbQ29udGVudF9U&#10;eXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhADj9If/WAAAAlAEAAAsAAAAAAAAAAAAAAAAALwEAAF9y&#10;ZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAFqKTa+gAgAAjgUAAA4AAAAAAAAAAAAAAAAALgIAAGRy&#10;cy9lMm9Eb2MueG1sUEsBAi0AFAAGAAgAAAAhAHIrnUrdAAAACQEAAA8AAAAAAAAAAAAAAAAA+gQA&#10;AGRycy9kb3ducmV2LnhtbFBLBQYAAAAABAAEAPMAAAAEBgAAAAA=&#10;" filled="f" strokecolor="black [3213]"/></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251667456" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="3F27F274" wp14:editId="7A8D2800"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>924560</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>287020</wp:posOffset></wp:positionV><wp:extent cx="297180" cy="228600"/><wp:effectExtent l="0" t="0" r="26670" b="19050"/><wp:wrapNone/><wp:docPr id="3" name="Rectángulo 3"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="297180" cy="228600"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:solidFill><a:schemeClr val="tx1"/></a:solidFill></a:ln></wps:spPr><wps:style><a:lnRef idx="2"><a:schemeClr val="accent1"><a:shade val="50000"/></a:schemeClr></a:lnRef><a:fillRef idx="1"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect w14:anchorId="1D653557" id="Rectángulo 3" o:spid="_x0000_s1026" style="position:absolute;margin-left:72.8pt;margin-top:22.6pt;width:23.4pt;height:18pt;z-index:251667456;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:middle" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQAaKnmgoQIAAI4FAAAOAAAAZHJzL2Uyb0RvYy54bWysVEtu2zAQ3RfoHQjuG30S52NEDowEKQoE&#10;aZCkyJqhSEsAxWFJ2rJ7m56lF+uQlGQjDboo6gXN0cy84bz5XF5tO0U2wroWdEWLo5wSoTnUrV5V&#10;9Nvz7adzSpxnumYKtKjoTjh6tfj44bI3c1FCA6oWliCIdvPeVLTx3syzzPFGdMwdgREalRJsxzyK&#10;dpXVlvWI3qmszPPTrAdbGwtcOIdfb5KSLiK+lIL7r1I64YmqKL7Nx9PG8zWc2eKSzVeWmablwzPY&#10;P7yiY63GoBPUDfOMrG37B1TXcgsOpD/i0GUgZctFzAGzKfI32Tw1zIiYC5LjzEST+3+w/H7zYElb&#10;V/SYEs06LNEjkvbrp16tFZDjQFBv3BztnsyDHSSH15DtVtou/GMeZBtJ3U2kiq0nHD+WF2fFOVLP&#10;UVWW56d5JD3bOxvr/GcBHQmXiloMH6lkmzvnMSCajiYhlobbVqlYN6VJX9GLWTmLDg5UWwdlMIsd&#10;JK6VJRuGtffbIqSCWAdWKCmNH0OCKaV48zslAoTSj0IiNyGJFCB05R6TcS60L5KqYbVIoWY5/sZg&#10;o0cMHQEDssRHTtgDwGiZQEbs9ObBPriK2NSTc/63hyXnySNGBu0n567VYN8DUJjVEDnZjyQlagJL&#10;r1DvsHMspJFyht+2WL875vwDszhDWHLcC/4rHlIB1gmGGyUN2B/vfQ/22NqopaTHmayo+75mVlCi&#10;vmhs+ovi5CQMcRROZmclCvZQ83qo0evuGrD0BW4gw+M12Hs1XqWF7gXXxzJERRXTHGNXlHs7Ctc+&#10;7QpcQFwsl9EMB9cwf6efDA/ggdXQn8/bF2bN0MQeu/8exvll8ze9nGyDp4bl2oNsY6PveR34xqGP&#10;jTMsqLBVDuVotV+ji98AAAD//wMAUEsDBBQABgAIAAAAIQApD2zV3QAAAAkBAAAPAAAAZHJzL2Rv&#10;d25yZXYueG1sTI9BTsMwEEX3SNzBmkpsEHVq0iqkcSqExDJItBzAjYc4ajx2Y6cNt8ddwfJrnv5/&#10;U+1mO7ALjqF3JGG1zIAhtU731En4Orw/FcBCVKTV4Agl/GCAXX1/V6lSuyt94mUfO5ZKKJRKgonR&#10;l5yH1qBVYek8Urp9u9GqmOLYcT2qayq3AxdZtuFW9ZQWjPL4ZrA97ScrYZ6K87mZTtbgczM8iug/&#10;Gu+lfFjMr1tgEef4B8NNP6lDnZyObiId2JByvt4kVEK+FsBuwIvIgR0lFCsBvK74/w/qXwAAAP//&#10;AwBQSwECLQAUAAYACAAAACEAtoM4kv4AAADhAQAAEwAAAAAAAAAAAAAAAAAAAAAAW0NvbnRlbnRf&#10;VHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAAAAAAAAAAAAAAC8BAABf&#10;cmVscy8ucmVsc1BLAQItABQABgAIAAAAIQAaKnmgoQIAAI4FAAAOAAAAAAAAAAAAAAAAAC4CAABk&#10;cnMvZTJvRG9jLnhtbFBLAQItABQABgAIAAAAIQApD2zV3QAAAAkBAAAPAAAAAAAAAAAAAAAAAPsE&#10;AABkcnMvZG93bnJldi54bWxQSwUGAAAAAAQABADzAAAABQYAAAAA&#10;" filled="f" strokecolor="black [3213]"/></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>Estado de la solicitud:</w:t></w:r></w:p><w:p w14:paraId="4D4C1AC3" w14:textId="2E9D014B" w:rsidR="007307BE" w:rsidRPr="00224C83" w:rsidRDefault="007307BE" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251666432" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="14447C5E" wp14:editId="2263D913"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>-635</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>-635</wp:posOffset></wp:positionV><wp:extent cx="297180" cy="228600"/><wp:effectExtent l="0" t="0" r="26670" b="19050"/><wp:wrapNone/><wp:docPr id="1" name="Rectángulo 1"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="297180" cy="228600"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:solidFill><a:schemeClr val="tx1"/></a:solidFill></a:ln></wps:spPr><wps:style><a:lnRef idx="2"><a:schemeClr val="accent1"><a:shade val="50000"/></a:schemeClr></a:lnRef><a:fillRef idx="1"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect w14:anchorId="74C826DE" id="Rectángulo 1" o:spid="_x0000_s1026" style="position:absolute;margin-left:-.05pt;margin-top:-.05pt;width:23.4pt;height:18pt;z-index:251666432;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:middle" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQBFnvVtoAIAAI4FAAAOAAAAZHJzL2Uyb0RvYy54bWysVEtu2zAQ3RfoHQjuG8mCnY8ROTASuCgQ&#10;JEE+yJqhSEsAxWFJ2rJ7m56lF+uQlGQjDboo6gXN0cy84bz5XF7tWkW2wroGdEknJzklQnOoGr0u&#10;6cvz6ss5Jc4zXTEFWpR0Lxy9Wnz+dNmZuSigBlUJSxBEu3lnSlp7b+ZZ5ngtWuZOwAiNSgm2ZR5F&#10;u84qyzpEb1VW5Plp1oGtjAUunMOvN0lJFxFfSsH9vZROeKJKim/z8bTxfAtntrhk87Vlpm54/wz2&#10;D69oWaMx6Ah1wzwjG9v8AdU23IID6U84tBlI2XARc8BsJvm7bJ5qZkTMBclxZqTJ/T9Yfrd9sKSp&#10;sHaUaNZiiR6RtF8/9XqjgEwCQZ1xc7R7Mg+2lxxeQ7Y7advwj3mQXSR1P5Iqdp5w/FhcnE3OkXqO&#10;qqI4P80j6dnB2VjnvwpoSbiU1GL4SCXb3jqPAdF0MAmxNKwapWLdlCZdSS9mxSw6OFBNFZTBLHaQ&#10;uFaWbBnW3u9iKoh1ZIWS0hggJJhSije/VyJAKP0oJHITkkgBQlceMBnnQvtJUtWsEinULMdf4C0E&#10;GzyiFAEDssRHjtg9wGCZQAbsBNPbB1cRm3p0zv/2sOQ8esTIoP3o3DYa7EcACrPqIyf7gaRETWDp&#10;Dao9do6FNFLO8FWD9btlzj8wizOEJce94O/xkAqwTtDfKKnB/vjoe7DH1kYtJR3OZEnd9w2zghL1&#10;TWPTX0ym0zDEUZjOzgoU7LHm7VijN+01YOmxsfF18RrsvRqu0kL7iutjGaKiimmOsUvKvR2Ea592&#10;BS4gLpbLaIaDa5i/1U+GB/DAaujP590rs6ZvYo/dfwfD/LL5u15OtsFTw3LjQTax0Q+89nzj0MfG&#10;6RdU2CrHcrQ6rNHFbwAAAP//AwBQSwMEFAAGAAgAAAAhANkWyNjZAAAABQEAAA8AAABkcnMvZG93&#10;bnJldi54bWxMjsFOwzAQRO9I/IO1SFxQ67SFUkKcCiFxDBKFD3DjJY5qr93YacPfs4gDnEajGc28&#10;ajt5J044pD6QgsW8AIHUBtNTp+Dj/WW2AZGyJqNdIFTwhQm29eVFpUsTzvSGp13uBI9QKrUCm3Ms&#10;pUytRa/TPEQkzj7D4HVmO3TSDPrM497JZVGspdc98YPVEZ8ttofd6BVM4+Z4bMaDt7hq3M0yx9cm&#10;RqWur6anRxAZp/xXhh98RoeamfZhJJOEUzBbcPFXOL1d34PYK1jdPYCsK/mfvv4GAAD//wMAUEsB&#10;Ai0AFAAGAAgAAAAhALaDOJL+AAAA4QEAABMAAAAAAAAAAAAAAAAAAAAAAFtDb250ZW50X1R5cGVz&#10;XS54bWxQSwECLQAUAAYACAAAACEAOP0h/9YAAACUAQAACwAAAAAAAAAAAAAAAAAvAQAAX3JlbHMv&#10;LnJlbHNQSwECLQAUAAYACAAAACEARZ71baACAACOBQAADgAAAAAAAAAAAAAAAAAuAgAAZHJzL2Uy&#10;b0RvYy54bWxQSwECLQAUAAYACAAAACEA2RbI2NkAAAAFAQAADwAAAAAAAAAAAAAAAAD6BAAAZHJz&#10;L2Rvd25yZXYueG1sUEsFBgAAAAAEAAQA8wAAAAAGAAAAAA==&#10;" filled="f" strokecolor="black [3213]"/></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r w:rsidRPr="00224C83"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">        </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Iniciada           Evaluada            </w:t></w:r><w:r w:rsidR="006A7529"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>Resuelta</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4662" w:type="dxa"/><w:gridSpan w:val="3"/><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="2136A223" w14:textId="260B9E56" w:rsidR="007307BE" w:rsidRPr="007307BE" w:rsidRDefault="007307BE" w:rsidP="007307BE"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>Fecha de la solicitud:</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00252BD9" w14:paraId="77B6A05F" w14:textId="77777777" w:rsidTr="00C8027A"><w:trPr><w:trHeight w:val="284"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="9746" w:type="dxa"/><w:gridSpan w:val="5"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/></w:tcPr><w:p w14:paraId="010D565B" w14:textId="77777777" w:rsidR="00252BD9" w:rsidRPr="00252BD9" w:rsidRDefault="00252BD9" w:rsidP="00FB5666"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t>Información de(los) solicitante(s)</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00D43AAC" w14:paraId="34885666" w14:textId="77777777" w:rsidTr="006A7529"><w:trPr><w:trHeight w:val="256"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5084" w:type="dxa"/><w:gridSpan w:val="2"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/></w:tcPr><w:p w14:paraId="34334903" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t>Nombre(s) y Apellido(s)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4662" w:type="dxa"/><w:gridSpan w:val="3"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/></w:tcPr><w:p w14:paraId="2E84B7B5" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t>Correo electrónico</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00A75C88" w14:paraId="00140102" w14:textId="77777777" w:rsidTr="006A7529"><w:trPr><w:trHeight w:val="256"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5084" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p w14:paraId="6D675102" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4662" w:type="dxa"/><w:gridSpan w:val="3"/></w:tcPr><w:p w14:paraId="2F9B9E48" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr w:rsidR="00A75C88" w14:paraId="23C81AE8" w14:textId="77777777" w:rsidTr="006A7529"><w:trPr><w:trHeight w:val="256"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5084" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p w14:paraId="1A775CD5" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4662" w:type="dxa"/><w:gridSpan w:val="3"/></w:tcPr><w:p w14:paraId="70310635" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr w:rsidR="00A75C88" w14:paraId="757107FB" w14:textId="77777777" w:rsidTr="006A7529"><w:trPr><w:trHeight w:val="256"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5084" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p w14:paraId="4EBCFD17" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4662" w:type="dxa"/><w:gridSpan w:val="3"/></w:tcPr><w:p w14:paraId="434627C8" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr w:rsidR="00A75C88" w14:paraId="18DC0336" w14:textId="77777777" w:rsidTr="006A7529"><w:trPr><w:trHeight w:val="256"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5084" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p w14:paraId="0C4C30B5" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4662" w:type="dxa"/><w:gridSpan w:val="3"/></w:tcPr><w:p w14:paraId="61F5416E" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr w:rsidR="00A75C88" w14:paraId="3568BE03" w14:textId="77777777" w:rsidTr="006A7529"><w:trPr><w:trHeight w:val="256"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5084" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p w14:paraId="4B73F7D2" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4662" w:type="dxa"/><w:gridSpan w:val="3"/></w:tcPr><w:p w14:paraId="3AF71EE3" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr w:rsidR="00A75C88" w14:paraId="19D5B018" w14:textId="77777777" w:rsidTr="006A7529"><w:trPr><w:trHeight w:val="256"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5084" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p w14:paraId="0E4E40E0" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4662" w:type="dxa"/><w:gridSpan w:val="3"/></w:tcPr><w:p w14:paraId="6B1B854C" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr w:rsidR="00A75C88" w14:paraId="444486E3" w14:textId="77777777" w:rsidTr="006A7529"><w:trPr><w:trHeight w:val="256"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5084" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p w14:paraId="35D14962" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4662" w:type="dxa"/><w:gridSpan w:val="3"/></w:tcPr><w:p w14:paraId="38F90513" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr w:rsidR="00A75C88" w14:paraId="3A8F66D4" w14:textId="77777777" w:rsidTr="006A7529"><w:trPr><w:trHeight w:val="256"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5084" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p w14:paraId="13F79C44" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4662" w:type="dxa"/><w:gridSpan w:val="3"/></w:tcPr><w:p w14:paraId="1D32A277" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr w:rsidR="00A75C88" w14:paraId="0CABCC4A" w14:textId="77777777" w:rsidTr="006A7529"><w:trPr><w:trHeight w:val="256"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5084" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p w14:paraId="2A68A451" w14:textId="77777777" w:rsidR="00D75519" w:rsidRDefault="00D75519" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4662" w:type="dxa"/><w:gridSpan w:val="3"/></w:tcPr><w:p w14:paraId="0B746DFA" w14:textId="77777777" w:rsidR="00A75C88" w:rsidRDefault="00A75C88" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr w:rsidR="00A75C88" w14:paraId="7A95A26C" w14:textId="77777777" w:rsidTr="00C8027A"><w:trPr><w:trHeight w:val="535"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="9746" w:type="dxa"/><w:gridSpan w:val="5"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/></w:tcPr><w:p w14:paraId="35C8FCC3" w14:textId="3C90D5F9" w:rsidR="00A75C88" w:rsidRDefault="00C8027A" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t>Caso(s) de uso</w:t></w:r><w:r w:rsidR="0086294F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> impactado(s):</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="0086294F" w14:paraId="5F074806" w14:textId="77777777" w:rsidTr="00C8027A"><w:trPr><w:trHeight w:val="914"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="9746" w:type="dxa"/><w:gridSpan w:val="5"/></w:tcPr><w:p w14:paraId="5985046A" w14:textId="77777777" w:rsidR="0086294F" w:rsidRDefault="0086294F" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr w:rsidR="0086294F" w14:paraId="2F7CED34" w14:textId="77777777" w:rsidTr="00C8027A"><w:trPr><w:trHeight w:val="497"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="9746" w:type="dxa"/><w:gridSpan w:val="5"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/></w:tcPr><w:p w14:paraId="3562EF23" w14:textId="3EC392F8" w:rsidR="0086294F" w:rsidRDefault="003D138E" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t>Nivel</w:t></w:r><w:r w:rsidR="00FB5666"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> del impacto en </w:t></w:r><w:r w:rsidR="00C8027A"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t>el(los) caso</w:t></w:r><w:r w:rsidR="00A522D9"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t>(s)</w:t></w:r><w:r w:rsidR="00C8027A"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> de uso existente(s)</w:t></w:r><w:r w:rsidR="00FB5666"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t>:</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00FB5666" w14:paraId="57D06EF6" w14:textId="77777777" w:rsidTr="00C8027A"><w:trPr><w:trHeight w:val="1020"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="9746" w:type="dxa"/><w:gridSpan w:val="5"/></w:tcPr><w:p w14:paraId="7CFBD4F5" w14:textId="6E4FC756" w:rsidR="00FB5666" w:rsidRPr="003D138E" w:rsidRDefault="003D138E" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251670528" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="74CED9D7" wp14:editId="7E202BD8"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>2951480</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>23495</wp:posOffset></wp:positionV><wp:extent cx="297180" cy="228600"/><wp:effectExtent l="0" t="0" r="26670" b="19050"/><wp:wrapNone/><wp:docPr id="2" name="Rectángulo 2"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="297180" cy="228600"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:solidFill><a:schemeClr val="tx1"/></a:solidFill></a:ln></wps:spPr><wps:style><a:lnRef idx="2"><a:schemeClr val="accent1"><a:shade val="50000"/></a:schemeClr></a:lnRef><a:fillRef idx="1"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect w14:anchorId="0210F3C4" id="Rectángulo 2" o:spid="_x0000_s1026" style="position:absolute;margin-left:232.4pt;margin-top:1.85pt;width:23.4pt;height:18pt;z-index:251670528;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:middle" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQAV8wcroAIAAI4FAAAOAAAAZHJzL2Uyb0RvYy54bWysVMlu2zAQvRfoPxC8N1oQZzEiB0aCFAWC&#10;JMiCnBmKtARQHJakLbt/02/pj3VISrKRBj0U9YHmaGbecN4sF5fbTpGNsK4FXdHiKKdEaA51q1cV&#10;fXm++XJGifNM10yBFhXdCUcvF58/XfRmLkpoQNXCEgTRbt6bijbem3mWOd6IjrkjMEKjUoLtmEfR&#10;rrLash7RO5WVeX6S9WBrY4EL5/DrdVLSRcSXUnB/L6UTnqiK4tt8PG0838KZLS7YfGWZaVo+PIP9&#10;wys61moMOkFdM8/I2rZ/QHUtt+BA+iMOXQZStlzEHDCbIn+XzVPDjIi5IDnOTDS5/wfL7zYPlrR1&#10;RUtKNOuwRI9I2q+ferVWQMpAUG/cHO2ezIMdJIfXkO1W2i78Yx5kG0ndTaSKrSccP5bnp8UZUs9R&#10;VZZnJ3kkPds7G+v8VwEdCZeKWgwfqWSbW+cxIJqOJiGWhptWqVg3pUlf0fNZOYsODlRbB2Uwix0k&#10;rpQlG4a199sipIJYB1YoKY0fQ4IppXjzOyUChNKPQiI3IYkUIHTlHpNxLrQvkqphtUihZjn+xmCj&#10;RwwdAQOyxEdO2APAaJlARuz05sE+uIrY1JNz/reHJefJI0YG7SfnrtVgPwJQmNUQOdmPJCVqAktv&#10;UO+wcyykkXKG37RYv1vm/AOzOENYctwL/h4PqQDrBMONkgbsj4++B3tsbdRS0uNMVtR9XzMrKFHf&#10;NDb9eXF8HIY4Csez0xIFe6h5O9TodXcFWPoCN5Dh8RrsvRqv0kL3iutjGaKiimmOsSvKvR2FK592&#10;BS4gLpbLaIaDa5i/1U+GB/DAaujP5+0rs2ZoYo/dfwfj/LL5u15OtsFTw3LtQbax0fe8Dnzj0MfG&#10;GRZU2CqHcrTar9HFbwAAAP//AwBQSwMEFAAGAAgAAAAhAD+4AlHcAAAACAEAAA8AAABkcnMvZG93&#10;bnJldi54bWxMj8FOwzAQRO9I/IO1SFwQddKWtIQ4FULiGCQKH+DGSxw1Xrux04a/ZznBcTSjmTfV&#10;bnaDOOMYe08K8kUGAqn1pqdOwefH6/0WREyajB48oYJvjLCrr68qXRp/oXc871MnuIRiqRXYlEIp&#10;ZWwtOh0XPiCx9+VHpxPLsZNm1Bcud4NcZlkhne6JF6wO+GKxPe4np2CetqdTMx2dxVUz3C1TeGtC&#10;UOr2Zn5+ApFwTn9h+MVndKiZ6eAnMlEMCtbFmtGTgtUGBPsPeV6AOLB+3ICsK/n/QP0DAAD//wMA&#10;UEsBAi0AFAAGAAgAAAAhALaDOJL+AAAA4QEAABMAAAAAAAAAAAAAAAAAAAAAAFtDb250ZW50X1R5&#10;cGVzXS54bWxQSwECLQAUAAYACAAAACEAOP0h/9YAAACUAQAACwAAAAAAAAAAAAAAAAAvAQAAX3Jl&#10;bHMvLnJlbHNQSwECLQAUAAYACAAAACEAFfMHK6ACAACOBQAADgAAAAAAAAAAAAAAAAAuAgAAZHJz&#10;L2Uyb0RvYy54bWxQSwECLQAUAAYACAAAACEAP7gCUdwAAAAIAQAADwAAAAAAAAAAAAAAAAD6BAAA&#10;ZHJzL2Rvd25yZXYueG1sUEsFBgAAAAAEAAQA8wAAAAMGAAAAAA==&#10;" filled="f" strokecolor="black [3213]"/></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251672576" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="05379DE4" wp14:editId="292F67A6"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>3332480</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>268605</wp:posOffset></wp:positionV><wp:extent cx="297180" cy="228600"/><wp:effectExtent l="0" t="0" r="26670" b="19050"/><wp:wrapNone/><wp:docPr id="4" name="Rectángulo 4"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="297180" cy="228600"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:solidFill><a:schemeClr val="tx1"/></a:solidFill></a:ln></wps:spPr><wps:style><a:lnRef idx="2"><a:schemeClr val="accent1"><a:shade val="50000"/></a:schemeClr></a:lnRef><a:fillRef idx="1"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect w14:anchorId="65410222" id="Rectángulo 4" o:spid="_x0000_s1026" style="position:absolute;margin-left:262.4pt;margin-top:21.15pt;width:23.4pt;height:18pt;z-index:251672576;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:middle" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQC1KeOmoAIAAI4FAAAOAAAAZHJzL2Uyb0RvYy54bWysVNtqGzEQfS/0H4Tem71g52KyDiYhpRDS&#10;kKTkWdFKXoFWo0qy1+7f9Fv6Yx1pLzZp6EOpH2TNzswZzZnL5dWu1WQrnFdgKlqc5JQIw6FWZl3R&#10;b8+3n84p8YGZmmkwoqJ74enV8uOHy84uRAkN6Fo4giDGLzpb0SYEu8gyzxvRMn8CVhhUSnAtCyi6&#10;dVY71iF6q7Myz0+zDlxtHXDhPX696ZV0mfClFDx8ldKLQHRF8W0hnS6dr/HMlpdssXbMNooPz2D/&#10;8IqWKYNBJ6gbFhjZOPUHVKu4Aw8ynHBoM5BScZFywGyK/E02Tw2zIuWC5Hg70eT/Hyy/3z44ouqK&#10;zigxrMUSPSJpv36a9UYDmUWCOusXaPdkH9wgebzGbHfStfEf8yC7ROp+IlXsAuH4sbw4K86Reo6q&#10;sjw/zRPp2cHZOh8+C2hJvFTUYfhEJdve+YAB0XQ0ibEM3CqtU920IV1FL+blPDl40KqOymiWOkhc&#10;a0e2DGsfdkVMBbGOrFDSBj/GBPuU0i3stYgQ2jwKidzEJPoAsSsPmIxzYULRqxpWiz7UPMffGGz0&#10;SKETYESW+MgJewAYLXuQEbt/82AfXUVq6sk5/9vDeufJI0UGEybnVhlw7wFozGqI3NuPJPXURJZe&#10;od5j5zjoR8pbfquwfnfMhwfmcIaw5LgXwlc8pAasEww3ShpwP977Hu2xtVFLSYczWVH/fcOcoER/&#10;Mdj0F8VsFoc4CbP5WYmCO9a8HmvMpr0GLH2BG8jydI32QY9X6aB9wfWxilFRxQzH2BXlwY3Cdeh3&#10;BS4gLlarZIaDa1m4M0+WR/DIauzP590Lc3Zo4oDdfw/j/LLFm17ubaOngdUmgFSp0Q+8Dnzj0KfG&#10;GRZU3CrHcrI6rNHlbwAAAP//AwBQSwMEFAAGAAgAAAAhAPDwmnbeAAAACQEAAA8AAABkcnMvZG93&#10;bnJldi54bWxMj81OwzAQhO9IvIO1lbgg6jTpT5TGqRASxyBReAA3XuKo9tqNnTa8PeYEx9GMZr6p&#10;D7M17IpjGBwJWC0zYEidUwP1Aj4/Xp9KYCFKUtI4QgHfGODQ3N/VslLuRu94PcaepRIKlRSgY/QV&#10;56HTaGVYOo+UvC83WhmTHHuuRnlL5dbwPMu23MqB0oKWHl80dufjZAXMU3m5tNPZaixa85hH/9Z6&#10;L8TDYn7eA4s4x78w/OIndGgS08lNpAIzAjb5OqFHAeu8AJYCm91qC+wkYFcWwJua/3/Q/AAAAP//&#10;AwBQSwECLQAUAAYACAAAACEAtoM4kv4AAADhAQAAEwAAAAAAAAAAAAAAAAAAAAAAW0NvbnRlbnRf&#10;VHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAAAAAAAAAAAAAAC8BAABf&#10;cmVscy8ucmVsc1BLAQItABQABgAIAAAAIQC1KeOmoAIAAI4FAAAOAAAAAAAAAAAAAAAAAC4CAABk&#10;cnMvZTJvRG9jLnhtbFBLAQItABQABgAIAAAAIQDw8Jp23gAAAAkBAAAPAAAAAAAAAAAAAAAAAPoE&#10;AABkcnMvZG93bnJldi54bWxQSwUGAAAAAAQABADzAAAABQYAAAAA&#10;" filled="f" strokecolor="black [3213]"/></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r w:rsidRPr="003D138E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t>Urgente (en un día debe estar resuelto máximo)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p w14:paraId="66A4B569" w14:textId="059F2C0D" w:rsidR="003D138E" w:rsidRPr="003D138E" w:rsidRDefault="003D138E" w:rsidP="00FB5666"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:noProof/><w:sz w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="003D138E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t>Mediano (se necesita estar resuelto máximo en 3 días)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:noProof/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00FB5666" w14:paraId="6400C88F" w14:textId="77777777" w:rsidTr="00C8027A"><w:tblPrEx><w:tblBorders><w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:insideH w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:insideV w:val="none" w:sz="0" w:space="0" w:color="auto"/></w:tblBorders></w:tblPrEx><w:trPr><w:trHeight w:val="51"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="9746" w:type="dxa"/><w:gridSpan w:val="5"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p w14:paraId="327F9841" w14:textId="77777777" w:rsidR="00FB5666" w:rsidRDefault="00FB5666" w:rsidP="00D43AAC"><w:pPr><w:tabs><w:tab w:val="left" w:pos="6636"/></w:tabs></w:pPr></w:p></w:tc></w:tr></w:tbl><w:p w14:paraId="0D1AF2B7" w14:textId="77777777" w:rsidR="00BA1882" w:rsidRPr="00BA1882" w:rsidRDefault="00BA1882" w:rsidP="0051501F"><w:pPr><w:tabs><w:tab w:val="left" w:pos="984"/></w:tabs></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)
